# Handback report generation: mark the two handed-off files as handed back,
# in sync with en-US, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$fileA = "6f0ac679-5966-4a62-884b-932e5d1582ef.md"
$fileB = "cde4fe71-6ed3-4d87-a2b9-13314eb40857.md"
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30a420171f83ca8cf3fe52deaa0c71d7fea5d516/e2e/6f0ac679-5966-4a62-884b-932e5d1582ef.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30a420171f83ca8cf3fe52deaa0c71d7fea5d516/e2e/cde4fe71-6ed3-4d87-a2b9-13314eb40857.md"
$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update the status column (E/F) for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status

# Widen the (now much longer) status columns to match.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Per-locale sheets (zh-cn / de-de): record the handback target file,
# handback xliff file (hyperlinked), and handback datetime.
# ---------------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackA = "6f0ac679-5966-4a62-884b-932e5d1582ef.b535a8a7665756ef3960379e473809fb14b83416.zh-cn.xlf"; HandbackB = "cde4fe71-6ed3-4d87-a2b9-13314eb40857.65965c95a06dfb0ce43fce7963b0a372177aefdc.zh-cn.xlf"; Datetime = "2016-08-27 08:24:37" },
    @{ Name = "de-de"; HandbackA = "6f0ac679-5966-4a62-884b-932e5d1582ef.b535a8a7665756ef3960379e473809fb14b83416.de-de.xlf"; HandbackB = "cde4fe71-6ed3-4d87-a2b9-13314eb40857.65965c95a06dfb0ce43fce7963b0a372177aefdc.de-de.xlf"; Datetime = "2016-08-27 08:24:44" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column (Status)
    $ws.Range("C2").Value = $status
    $ws.Range("C3").Value = $status

    # Latest Target File (I) -- the handed-back source file, hyperlinked.
    $ws.Range("I2").Value = $fileA
    $ws.Range("I3").Value = $fileB
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, "", "", $fileA)
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlB, "", "", $fileB)

    # Latest Handback File (J)
    $ws.Range("J2").Value = $locale.HandbackA
    $ws.Range("J3").Value = $locale.HandbackB

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $locale.Datetime
    $ws.Range("K3").Value = $locale.Datetime

    # Widen Status (C) and Latest Target/Handback File (I/J) columns.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
